$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Mace/*Warhammer=Light Crossbow!20 bolts/?SW"
$ws.Range("A6").Value = "1=Divine Domain:=classes/cleric/DivineDomains.xlsx"
$ws.Range("B6").Value = "3=Roguish Archetype:=classes/rogue/RoguishArchetypes.xlsx"
